$wb = $excel.ActiveWorkbook

function Set-CellValues {
    param($ws, $pairs)
    foreach ($pair in $pairs) {
        $cellRef = $pair[0]
        $val = $pair[1]
        if ($val -eq $null) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $val
        }
    }
}

$ws = $wb.Worksheets.Item("ALC")
$ALC_data = @(
    @("H62", 10827.448),
    @("I62", 15793.066),
    @("J62", 5507.143),
    @("K62", 15793.066),
    @("L62", 5507.143),
    @("M62", -15169.066),
    @("N62", -6755.143),
    @("H64", 4426.7036),
    @("I64", 4181.636),
    @("J64", 4595.1875),
    @("K64", 4181.636),
    @("L64", 4595.1875),
    @("M64", -3933.636),
    @("N64", -5091.1875),
    @("H65", 10827.448),
    @("I65", 15793.066),
    @("J65", 5507.143),
    @("K65", 78965.33),
    @("L65", 27535.715),
    @("M65", -75845.33),
    @("N65", -33775.715),
    @("H67", 4426.7036),
    @("I67", 4181.636),
    @("J67", 4595.1875),
    @("K67", 4181.636),
    @("L67", 4595.1875),
    @("M67", -3323.636),
    @("N67", -6311.1875),
    @("H76", 20841516),
    @("I76", 62515400),
    @("J76", 4573.75),
    @("K76", 62515400),
    @("L76", 4573.75),
    @("M76", -62515085),
    @("N76", -5203.75),
    @("H79", 20841516),
    @("I79", 62515400),
    @("J79", 4573.75),
    @("K79", 62515400),
    @("L79", 4573.75),
    @("M79", -62514308),
    @("N79", -6757.75),
    @("H82", 602.2),
    @("I82", 602.2),
    @("K82", 1806.6),
    @("M82", -1400.6),
    @("H85", 602.2),
    @("I85", 602.2),
    @("K85", 1806.6),
    @("M85", -402.6000000000001),
    @("H116", 3179.3684),
    @("I116", 3128.25),
    @("K116", 3128.25),
    @("M116", 313.75),
    @("H132", 4606.74),
    @("I132", 1666.7561),
    @("K132", 5000.2683),
    @("M132", -2470.2683),
    @("H138", 1921.0878),
    @("I138", 1297.8334),
    @("J138", 2613.5925),
    @("K138", 3893.5002),
    @("L138", 7840.7775),
    @("M138", 1246.4998),
    @("N138", -18120.7775)
)
Set-CellValues -ws $ws -pairs $ALC_data

$ws = $wb.Worksheets.Item("ARM")
$ARM_data = @(
    @("H2", 2126.25),
    @("I2", 2819.4),
    @("K2", 2819.4),
    @("M2", -2706.4),
    @("H32", 1305338.2),
    @("I32", 1415644.9),
    @("J32", 3720),
    @("K32", 1415644.9),
    @("L32", 3720),
    @("M32", -1415357.9),
    @("N32", -4294),
    @("H45", 1000),
    @("I45", 0),
    @("J45", 1000),
    @("K45", 0),
    @("L45", 1000),
    @("M45", $null),
    @("N45", -1754),
    @("H61", 411337.5),
    @("I61", 295611.34),
    @("J61", 673650.0600000001),
    @("K61", 295611.34),
    @("L61", 673650.0600000001),
    @("M61", -295399.34),
    @("N61", -674074.0600000001),
    @("H116", 2126.25),
    @("I116", 2819.4),
    @("K116", 2819.4),
    @("M116", -525.4000000000001),
    @("H122", 3122.5),
    @("I122", 2069.2727),
    @("J122", 6984.3335),
    @("K122", 6207.8181),
    @("L122", 20953.0005),
    @("M122", -3757.8181),
    @("N122", -25853.0005),
    @("H136", 411337.5),
    @("I136", 295611.34),
    @("J136", 673650.0600000001),
    @("K136", 886834.02),
    @("L136", 2020950.18),
    @("M136", -884284.02),
    @("N136", -2026050.18)
)
Set-CellValues -ws $ws -pairs $ARM_data

$ws = $wb.Worksheets.Item("BSM")
$BSM_data = @(
    @("H3", 2126.25),
    @("I3", 2819.4),
    @("K3", 2819.4),
    @("M3", -2705.4),
    @("H80", 274.22223),
    @("I80", 339.16666),
    @("J80", 144.33333),
    @("K80", 339.16666),
    @("L80", 144.33333),
    @("M80", 658.83334),
    @("N80", -2140.33333),
    @("H83", 274.22223),
    @("I83", 339.16666),
    @("J83", 144.33333),
    @("K83", 1695.8333),
    @("L83", 721.6666499999999),
    @("M83", 3296.1667),
    @("N83", -10705.66665),
    @("H96", 29000),
    @("I96", 0),
    @("K96", 0),
    @("M96", $null),
    @("H99", 9901232),
    @("I99", 3510549.8),
    @("J99", 33333734),
    @("K99", 3510549.8),
    @("L99", 33333734),
    @("M99", -3509051.8),
    @("N99", -33336730),
    @("H107", 774.41174),
    @("I107", 774.41174),
    @("J107", 0),
    @("K107", 774.41174),
    @("L107", 0),
    @("M107", 1145.58826),
    @("N107", $null)
)
Set-CellValues -ws $ws -pairs $BSM_data

$ws = $wb.Worksheets.Item("CRP")
$CRP_data = @(
    @("H31", 920181.25),
    @("I31", 693.8982999999999),
    @("K31", 693.8982999999999),
    @("M31", -398.8982999999999),
    @("H34", 920181.25),
    @("I34", 693.8982999999999),
    @("K34", 693.8982999999999),
    @("M34", -491.8982999999999),
    @("H86", 4655.2173),
    @("I86", 2833.6428),
    @("J86", 7488.778),
    @("K86", 2833.6428),
    @("L86", 7488.778),
    @("M86", -1710.6428),
    @("N86", -9734.778),
    @("H89", 4655.2173),
    @("I89", 2833.6428),
    @("J89", 7488.778),
    @("K89", 14168.214),
    @("L89", 37443.89),
    @("M89", -8552.214),
    @("N89", -48675.89),
    @("H94", 10531),
    @("I94", 1835.8),
    @("K94", 1835.8),
    @("M94", -1384.8),
    @("H105", 1136.7028),
    @("I105", 980.72),
    @("K105", 980.72),
    @("M105", 766.28),
    @("H134", 11906644),
    @("I134", 15626457),
    @("J134", 3240.2),
    @("K134", 46879371),
    @("L134", 9720.599999999999),
    @("M134", -46876836),
    @("N134", -14790.6)
)
Set-CellValues -ws $ws -pairs $CRP_data

$ws = $wb.Worksheets.Item("CUL")
$CUL_data = @(
    @("H20", 0),
    @("J20", 0),
    @("L20", 0),
    @("N20", $null),
    @("H68", 1122.2903),
    @("I68", 731.2778),
    @("J68", 1282.25),
    @("K68", 2193.8334),
    @("L68", 3846.75),
    @("M68", -1382.8334),
    @("N68", -5468.75),
    @("H71", 1122.2903),
    @("I71", 731.2778),
    @("J71", 1282.25),
    @("K71", 6581.500199999999),
    @("L71", 11540.25),
    @("M71", -2525.500199999999),
    @("N71", -19652.25),
    @("H80", 1952.8572),
    @("J80", 1952.8572),
    @("L80", 5858.571599999999),
    @("N80", -7730.571599999999),
    @("H83", 1952.8572),
    @("J83", 1952.8572),
    @("L83", 17575.7148),
    @("N83", -26935.7148),
    @("H107", 1067.3704),
    @("I107", 288),
    @("J107", 2041.5834),
    @("K107", 864),
    @("L107", 6124.7502),
    @("M107", 1056),
    @("N107", -9964.7502),
    @("H132", 3233.3333),
    @("I132", 4980),
    @("J132", 1985.7142),
    @("K132", 44820),
    @("L132", 17871.4278),
    @("M132", -42290),
    @("N132", -22931.4278)
)
Set-CellValues -ws $ws -pairs $CUL_data

$ws = $wb.Worksheets.Item("GSM")
$GSM_data = @(
    @("H122", 1102.8462),
    @("I122", 1179.625),
    @("J122", 980),
    @("K122", 3538.875),
    @("L122", 2940),
    @("M122", -1088.875),
    @("N122", -7840),
    @("H126", 2469),
    @("I126", 1876.3636),
    @("J126", 3400.2856),
    @("K126", 5629.0908),
    @("L126", 10200.8568),
    @("M126", -3159.0908),
    @("N126", -15140.8568)
)
Set-CellValues -ws $ws -pairs $GSM_data

$ws = $wb.Worksheets.Item("WVR")
$WVR_data = @(
    @("H62", 4125),
    @("I62", 2333.3333),
    @("J62", 4538.4614),
    @("K62", 2333.3333),
    @("L62", 4538.4614),
    @("M62", -1709.3333),
    @("N62", -5786.4614),
    @("H65", 4125),
    @("I65", 2333.3333),
    @("J65", 4538.4614),
    @("K65", 11666.6665),
    @("L65", 22692.307),
    @("M65", -8546.666499999999),
    @("N65", -28932.307),
    @("H136", 20781540),
    @("I136", 31598620),
    @("J136", 5047604.5),
    @("K136", 94795860),
    @("L136", 15142813.5),
    @("M136", -94793310),
    @("N136", -15147913.5)
)
Set-CellValues -ws $ws -pairs $WVR_data

Write-Host "Applied all Sheets updates via scheduled runner"